$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23, shifting existing rows 23-45 down to 24-46
$ws.Rows("23:23").Insert()

# Populate the newly inserted row 23 with the new data entry
$ws.Range("A23").Value = 8
$ws.Range("B23").Value = "Terminal La Palmera de La Serena"
$ws.Range("C23").Value = "Coquimbo"
$ws.Range("D23").Value = 44740
$ws.Range("E23").Value = 4
$ws.Range("F23").Value = 100114007
$ws.Range("G23").Value = "Jengibre"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 500
$ws.Range("K23").Value = 16000
$ws.Range("L23").Value = 17000
$ws.Range("M23").Value = 16500
$ws.Range("N23").Value = "$/caja 13 kilos"
$ws.Range("O23").Value = "Perú"
$ws.Range("P23").Value = 1269
$ws.Range("Q23").Value = 13
$ws.Range("R23").Value = "Hortaliza"

# Match the date cell style used by the rest of column D
$ws.Range("D24").Copy()
$ws.Range("D23").PasteSpecial(-4122)
